$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Local Test" score value (F3); dependent formulas F4 (F2-F3) and
# F5 (F4/F2) will recalculate automatically.
$ws.Range("F3").Value = 1522529

# Update the active selection to match the new cursor position.
$ws.Range("H8").Select()
